$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = '[Achraf%Ammar%NULL%1,   Hamdi%Chtourou%NULL%2,   Hamdi%Chtourou%NULL%0,   Omar%Boukhris%NULL%2,   Omar%Boukhris%NULL%0,   Khaled%Trabelsi%NULL%2,   Khaled%Trabelsi%NULL%0,   Liwa%Masmoudi%NULL%2,   Liwa%Masmoudi%NULL%0,   Michael%Brach%NULL%1,   Bassem%Bouaziz%NULL%1,   Ellen%Bentlage%NULL%1,   Daniella%How%NULL%2,   Daniella%How%NULL%0,   Mona%Ahmed%NULL%1,   Patrick%Mueller%NULL%1,   Notger%Mueller%NULL%1,   Hsen%Hsouna%NULL%1,   Asma%Aloui%NULL%1,   Omar%Hammouda%NULL%1,   Laisa Liane%Paineiras-Domingos%NULL%1,   Annemarie%Braakman-Jansen%NULL%2,   Annemarie%Braakman-Jansen%NULL%0,   Christian%Wrede%NULL%1,   Sophia%Bastoni%NULL%1,   Carlos Soares%Pernambuco%NULL%1,   Leonardo Jose%Mataruna-Dos-Santos%NULL%1,   Morteza%Taheri%NULL%2,   Morteza%Taheri%NULL%0,   Khadijeh%Irandoust%NULL%2,   Khadijeh%Irandoust%NULL%0,   Aïmen%Khacharem%NULL%1,   Nicola L.%Bragazzi%NULL%1,   Jana%Strahler%NULL%2,   Jana%Strahler%NULL%0,   Jad Adrian%Washif%NULL%1,   Albina%Andreeva%NULL%1,   Samira C.%khoshnami%NULL%1,   Evangelia%Samara%NULL%1,   Vasiliki%Zisi%NULL%1,   Parasanth%Sankar%NULL%1,   Waseem N.%Ahmed%NULL%1,   Mohamed%Romdhani%NULL%1,   Jan%Delhey%NULL%1,   Stephen J.%Bailey%NULL%1,   Nicholas T.%Bott%NULL%1,   Faiez%Gargouri%NULL%2,   Faiez%Gargouri%NULL%0,   Lotfi%Chaari%NULL%1,   Hadj%Batatia%NULL%1,   Gamal Mohamed%Ali%NULL%1,   Osama%Abdelkarim%NULL%1,   Mohamed%Jarraya%NULL%1,   Kais El%Abed%NULL%1,   Nizar%Souissi%NULL%1,   Lisette Van%Gemert-Pijnen%NULL%1,   Bryan L.%Riemann%NULL%1,   Laurel%Riemann%NULL%1,   Wassim%Moalla%NULL%1,   Jonathan%Gómez-Raja%NULL%1,   Monique%Epstein%NULL%2,   Monique%Epstein%NULL%0,   Robbert%Sanderman%NULL%1,   Sebastian%Schulz%NULL%1,   Achim%Jerg%NULL%2,   Achim%Jerg%NULL%0,   Ramzi%Al-Horani%NULL%1,   Taiysir%Mansi%NULL%1,   Mohamed%Jmail%NULL%1,   Fernando%Barbosa%NULL%1,   Fernando%Ferreira-Santos%NULL%1,   Boštjan%Šimunič%NULL%2,   Boštjan%Šimunič%NULL%0,   Rado%Pišot%NULL%1,   Saša%Pišot%NULL%1,   Andrea%Gaggioli%NULL%1,   Piotr%Zmijewski%NULL%1,   Christian%Apfelbacher%NULL%1,   Jürgen%Steinacker%NULL%1,   Helmi Ben%Saad%NULL%1,   Jordan M.%Glenn%NULL%1,   Karim%Chamari%NULL%1,   Tarak%Driss%NULL%1,   Anita%Hoekelmann%NULL%1,   NULL%NULL%NULL%0]'

$ws.Range("E5").Value = '[Luigi%Barrea%luigi.barrea@unina.it%1,   Gabriella%Pugliese%NULL%2,   Gabriella%Pugliese%NULL%0,   Lydia%Framondi%NULL%1,   Rossana%Di Matteo%NULL%1,   Daniela%Laudisio%NULL%1,   Silvia%Savastano%NULL%1,   Annamaria%Colao%NULL%1,   Giovanna%Muscogiuri%NULL%1]'

$ws.Range("I6").Value = '_CROSSREF'

$ws.Range("E7").Value = '[Gemma%Biviá-Roig%NULL%1,   Valentina Lucia%La Rosa%NULL%2,   Valentina Lucia%La Rosa%NULL%0,   María%Gómez-Tébar%NULL%2,   María%Gómez-Tébar%NULL%0,   Lola%Serrano-Raya%NULL%1,   Juan José%Amer-Cuenca%NULL%1,   Salvatore%Caruso%NULL%1,   Elena%Commodari%NULL%2,   Elena%Commodari%NULL%0,   Antonio%Barrasa-Shaw%NULL%1,   Juan Francisco%Lisón%NULL%1]'

$ws.Range("E8").Value = '[Dimitrios I.%Bourdas%NULL%2,   Emmanouil D.%Zacharakis%NULL%3,   Emmanouil D.%Zacharakis%NULL%0]'

$ws.Range("E9").Value = '[Dimitrios I.%Bourdas%dbourdas@hotmail.gr%0,   Emmanouil D.%Zacharakis%NULL%0]'

$ws.Range("E10").Value = '[Collins%C.%coreGivesNoEmail%1,   Pavlidis%A.%coreGivesNoEmail%1,   Power%K.%coreGivesNoEmail%1,   Rowe%D.%coreGivesNoEmail%1,   van%Dalen H. P.%coreGivesNoEmail%1]'

$ws.Range("E11").Value = '[Dawn%Branley-Bell%dawn.branley-bell@northumbria.ac.uk%1,   Catherine V.%Talbot%NULL%2,   Catherine V.%Talbot%NULL%0]'

$ws.Range("E12").Value = '[Daniel D.%Callow%NULL%1,   Naomi A.%Arnold-Nedimala%NULL%1,   Leslie S.%Jordan%NULL%1,   Gabriel S.%Pena%NULL%1,   Junyeon%Won%NULL%1,   John L.%Woodard%NULL%1,   J. Carson%Smith%NULL%1]'

$ws.Range("E13").Value = '[Raffaella%Cancello%NULL%1,   Davide%Soranna%NULL%2,   Davide%Soranna%NULL%0,   Gaia%Zambra%NULL%1,   Antonella%Zambon%NULL%1,   Cecilia%Invitti%NULL%1]'

$ws.Range("E14").Value = '[Caruso%I.%coreGivesNoEmail%0,   Cignarelli%A.%coreGivesNoEmail%1,   Di%Molfetta S.%coreGivesNoEmail%1,   Giordano%F.%coreGivesNoEmail%1,   Giorgino%F.%coreGivesNoEmail%1,   Guarini%F.%coreGivesNoEmail%1,   Laviola%L.%coreGivesNoEmail%1,   Leonardini%A.%coreGivesNoEmail%1,   Natalicchio%A.%coreGivesNoEmail%1,   Perrini%S.%coreGivesNoEmail%1]'

$ws.Range("E15").Value = '[Arkaitz%Castañeda-Babarro%NULL%1,   Ane%Arbillaga-Etxarri%NULL%2,   Ane%Arbillaga-Etxarri%NULL%0,   Borja%Gutiérrez-Santamaría%NULL%1,   Aitor%Coca%NULL%2,   Aitor%Coca%NULL%0]'
